$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value from 45558 to 45559 for existing
# data rows 2-29.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Value = 45559
}

# Row 29 picks up an explicit row height (15, custom height) in the new file.
$ws.Rows.Item(29).RowHeight = 15

# Append the new data row (row 30).
$ws.Cells.Item(30, 1).Value = "A 40747-2024"
$ws.Cells.Item(30, 2).Value = 45558
$ws.Cells.Item(30, 3).Value = 45559
$ws.Cells.Item(30, 4).Value = "OKÄNT"
$ws.Cells.Item(30, 5).Value = "OKÄNT"
$ws.Cells.Item(30, 7).Value = 0.7
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = 0
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0

# Match formatting used by the rest of the table:
#  - column B/C use the date style (same as row 29)
#  - column R keeps the wrap-text style used throughout, left blank
$ws.Range("B30:C30").NumberFormat = $ws.Range("B29:C29").NumberFormat
$ws.Range("R30").WrapText = $true
